$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has a bold/centered/bordered header row (style 1),
# a bold wrap-text corner cell (style 2), and wrap-text cells down column G
# (style 3), plus an explicit row height of 32 on every row. The edit clears
# all of that direct formatting back to the workbook's default "Normal"
# style and lets row heights return to the sheet's standard height.

$dataRange = $ws.Range("A1:G19")

# Reset every cell back to the default "Normal" cell style - this removes
# the bold fonts, the thin borders and the wrap/centered alignment in one
# shot without disturbing any cell values.
$dataRange.Style = "Normal"

# With the wrap-text/formatting gone, auto-fit the rows so the explicit
# ht="32" row heights collapse back down to the sheet's standard height.
$dataRange.EntireRow.AutoFit()

# Move the active selection to K14, matching the new selection recorded on
# the sheet.
$ws.Range("K14").Select()
